$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.466
$ws.Range("C2").Value = 12.08
$ws.Range("D2").Value = 28.85
$ws.Range("E2").Value = 0.7682283586720448
$ws.Range("F2").Value = 0.2511848341232228

$ws.Range("E3").Value = 0.5442134657151719

$ws.Range("B4").Value = 12.514
$ws.Range("G4").Value = 0.4064773527340999

$ws.Range("E5").Value = 0.7071051815079117
$ws.Range("G5").Value = 0.7873738559023705

$ws.Range("G6").Value = 0.5817883126026757

$ws.Range("E7").Value = 0.9959664908470369
$ws.Range("G7").Value = 0.4832199014315892
